{"js": "// Locate the paragraph that ends with \"...research we could do to find\n// more concrete links between them.\" and:\n//   1. append a new run/sentence to it,\n//   2. insert a new blank paragraph after it,\n//   3. insert a \"Challenges and Limitations\" heading-like paragraph,\n//   4. insert the final \"biggest limitation...\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"research we could do to find more concrete links between them.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(marker) !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find target paragraph containing marker text.\");\n}\n\n// 1) Add the new sentence as an additional run at the end of the paragraph.\ntarget.insertText(\n  \" As it stands, the minor correlations between specific pollutants, and AQI\\u2019s of regions cannot confirm our hypothesis.\",\n  Word.InsertLocation.end\n);\n\n// 2) New empty paragraph directly after it.\nconst blankPara = target.insertParagraph(\"\", Word.InsertLocation.after);\n\n// 3) \"Challenges and Limitations\" paragraph after the blank one.\nconst headingPara = blankPara.insertParagraph(\"Challenges and Limitations\", Word.InsertLocation.after);\n\n// 4) Final paragraph with the limitations text.\nheadingPara.insertParagraph(\n  \"The biggest limitation we had was not collecting historical data. That was due to the cost involved in purchasing said data. Ideally, we would be able to collect this historic data and how that compares to the live data, and if there is any significantly noteworthy changes that events such as the COVID pandemic had between the data sets.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$rsq = [char]0x2019\n\n# Locate the paragraph ending in \"...research we could do to find more\n# concrete links between them.\" and append a new sentence/run to it.\n$range = $d.Content\n$found = $range.Find.Execute(\"research we could do to find more concrete links between them.\")\n\nif ($found) {\n    $range.Collapse(0)\n    $range.InsertAfter(\" As it stands, the minor correlations between specific pollutants, and AQI\" + $rsq + \"s of regions cannot confirm our hypothesis.\")\n\n    # New blank paragraph right after it.\n    $d.Paragraphs.Last.Range.InsertParagraphAfter()\n\n    # \"Challenges and Limitations\" paragraph.\n    $d.Paragraphs.Last.Range.InsertParagraphAfter()\n    $d.Paragraphs.Last.Range.Text = \"Challenges and Limitations\"\n\n    # Final limitations paragraph.\n    $d.Paragraphs.Last.Range.InsertParagraphAfter()\n    $d.Paragraphs.Last.Range.Text = \"The biggest limitation we had was not collecting historical data. That was due to the cost involved in purchasing said data. Ideally, we would be able to collect this historic data and how that compares to the live data, and if there is any significantly noteworthy changes that events such as the COVID pandemic had between the data sets.\"\n}\n"}
